$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a number by Excel;
# force them to Text format, assign, then restore the original (default) style
# so no numeric conversion happens and no stray formatting is left behind.
$textForceCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D10',
    'D11',
    'D15',
    'D18',
    'D20',
    'D21',
    'D22',
    'D23',
    'D28',
    'D29',
    'D31',
    'D32',
    'D33',
    'D38',
    'D39',
    'D41',
    'D42',
    'D43',
    'D47',
    'D48'
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '25.265.05'
$ws.Range('E2').Value = '  -2.88%  '
$ws.Range('D3').Value = '1.551.49'
$ws.Range('E3').Value = '  -4.82%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '206.90'
$ws.Range('E5').Value = '  -3.43%  '
$ws.Range('D6').Value = '1.00'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '0.478'
$ws.Range('E7').Value = '  -5.10%  '
$ws.Range('D8').Value = '0.0609'
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = '17.65'
$ws.Range('E10').Value = '  -4.89%  '
$ws.Range('D11').Value = '0.0780'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').Value = '1.765.82'
$ws.Range('E12').Value = '  -4.91%  '
$ws.Range('D13').Value = '1.550.57'
$ws.Range('E13').Value = '  -5.40%  '
$ws.Range('E14').Value = '  -4.71%  '
$ws.Range('D15').Value = '0.504'
$ws.Range('E15').Value = '  -4.53%  '
$ws.Range('D16').Value = '25.241.85'
$ws.Range('E16').Value = '  -2.99%  '
$ws.Range('D17').Value = '0.0₃0707'
$ws.Range('E17').Value = '  -4.80%  '
$ws.Range('D18').Value = '58.57'
$ws.Range('E18').Value = '  -4.92%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '185.45'
$ws.Range('E20').Value = '  -3.85%  '
$ws.Range('D21').Value = '4.10'
$ws.Range('E21').Value = '  -3.73%  '
$ws.Range('D22').Value = '9.25'
$ws.Range('E22').Value = '  -3.18%  '
$ws.Range('D23').Value = '5.84'
$ws.Range('E23').Value = '  -3.86%  '
$ws.Range('E24').Value = '  -4.37%  '
$ws.Range('E27').Value = '  -5.49%  '
$ws.Range('D28').Value = '14.84'
$ws.Range('E28').Value = '  -2.88%  '
$ws.Range('D29').Value = '6.38'
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('E30').Value = '  -7.08%  '
$ws.Range('D31').Value = '0.0464'
$ws.Range('E31').Value = '  -3.83%  '
$ws.Range('D32').Value = '3.02'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').Value = '2.97'
$ws.Range('E33').Value = '  -4.97%  '
$ws.Range('E34').Value = '  -3.16%  '
$ws.Range('E35').Value = '  -3.58%  '
$ws.Range('D36').Value = '1.081.52'
$ws.Range('E36').Value = '  -3.49%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = '0.0149'
$ws.Range('E38').Value = '  -3.08%  '
$ws.Range('D39').Value = '0.493'
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('E40').Value = '  -7.88%  '
$ws.Range('D41').Value = '0.762'
$ws.Range('E41').Value = '  -10.57%  '
$ws.Range('D42').Value = '0.797'
$ws.Range('E42').Value = '  +3.39%  '
$ws.Range('D43').Value = '92.43'
$ws.Range('E43').Value = '  -6.03%  '
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('D45').Value = '1.680.52'
$ws.Range('E45').Value = '  -4.88%  '
$ws.Range('E46').Value = '  -2.52%  '
$ws.Range('D47').Value = '1.46'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').Value = '52.24'
$ws.Range('E48').Value = '  -4.27%  '
$ws.Range('E49').Value = '  -5.06%  '
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('E51').Value = '  -2.41%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
